$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "DISTRITO"
$ws.Range("B1").Value = "RENTA_MEDIA"
$ws.Range("B2").Select()
